$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 20; this shifts the existing rows 20-32 down to 21-33
$ws.Rows("20:20").Insert()

# Populate the newly inserted row 20 with the new weekly data point.
# Values for columns that stay constant across all rows are copied from row 21
# (the row that used to be row 20 before the insert).
$ws.Range("A20").Value = 11
$ws.Range("B20").Value = "Vega Monumental Concepción"
$ws.Range("C20").Value = "Bíobío"
$ws.Range("D20").Value2 = 44741
$ws.Range("E20").Value = 8
$ws.Range("F20").Value = 100114007
$ws.Range("G20").Value = "Jengibre"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 50
$ws.Range("K20").Value = 14000
$ws.Range("L20").Value = 15000
$ws.Range("M20").Value = 14400
$ws.Range("N20").Value = "$/caja 13 kilos"
$ws.Range("O20").Value = "Perú"
$ws.Range("P20").Value = 1108
$ws.Range("Q20").Value = 13
$ws.Range("R20").Value = "Hortaliza"
